$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The existing rows 637:641 (last 5 data rows) are being "pushed down" to
# 642:646 unchanged, while rows 637:641 get refreshed with new weekly
# price figures. First duplicate the current rows 637:641 into the new
# rows 642:646 so the old data is preserved further down the sheet.
$src = $ws.Range("A637:T641")
$dst = $ws.Range("A642:T646")
$src.Copy($dst)

# Now overwrite rows 637:641 in place with the updated values.

# Row 637: Packham's Triumph, Calibre 70
$ws.Range("D637").Value = 44656
$ws.Range("L637").Value = "Calibre 70"
$ws.Range("M637").Value = 80
$ws.Range("N637").Value = 14000
$ws.Range("O637").Value = 14000
$ws.Range("P637").Value = 14000
$ws.Range("Q637").Value = '$/caja 18 kilos embalada'
$ws.Range("S637").Value = 778

# Row 638: Packham's Triumph, Calibre 90
$ws.Range("D638").Value = 44656
$ws.Range("L638").Value = "Calibre 90"
$ws.Range("M638").Value = 150
$ws.Range("N638").Value = 10000
$ws.Range("O638").Value = 10000
$ws.Range("P638").Value = 10000
$ws.Range("Q638").Value = '$/caja 18 kilos embalada'
$ws.Range("S638").Value = 556

# Row 639: Packham's Triumph, Primera (volume/price refreshed)
$ws.Range("D639").Value = 44656
$ws.Range("M639").Value = 200
$ws.Range("N639").Value = 12000
$ws.Range("O639").Value = 12000
$ws.Range("P639").Value = 12000
$ws.Range("S639").Value = 667

# Row 640: Winter Nelis, Especial
$ws.Range("D640").Value = 44656
$ws.Range("K640").Value = "Winter Nelis"
$ws.Range("L640").Value = "Especial"
$ws.Range("M640").Value = 80
$ws.Range("N640").Value = 18000
$ws.Range("O640").Value = 18000
$ws.Range("P640").Value = 18000
$ws.Range("Q640").Value = '$/bandeja 18 kilos granel'
$ws.Range("S640").Value = 1000
$ws.Range("T640").Value = 18

# Row 641: Winter Nelis, Primera
$ws.Range("D641").Value = 44656
$ws.Range("K641").Value = "Winter Nelis"
$ws.Range("L641").Value = "Primera"
$ws.Range("M641").Value = 150

$wb.Save()
